$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 39 (shifts rows 39:146 down to 40:147)
$ws.Rows(39).Insert()

# Populate the new row 39 with the new data record
$ws.Range("A39").Value = 11
$ws.Range("B39").Value = 'Vega Monumental Concepción'
$ws.Range("C39").Value = 'Bíobío'
$ws.Range("D39").Value = 44624
$ws.Range("E39").Value = 8
$ws.Range("F39").Value = 100112003
$ws.Range("G39").Value = 'Ajo'
$ws.Range("H39").Value = 'Chino'
$ws.Range("I39").Value = 'Primera'
$ws.Range("J39").Value = 220
$ws.Range("K39").Value = 16000
$ws.Range("L39").Value = 17000
$ws.Range("M39").Value = 16545
$ws.Range("N39").Value = '$/caja 10 kilos'
$ws.Range("O39").Value = 'China'
$ws.Range("P39").Value = 1654
$ws.Range("Q39").Value = 10
$ws.Range("R39").Value = 'Hortaliza'
